$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.910.33"
$ws.Range("E2").Value = "'  -1.10%  "
$ws.Range("D3").Value = "'1.634.89"
$ws.Range("E3").Value = "'  -2.59%  "
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("D5").Value = "'209.48"
$ws.Range("E5").Value = "'  -0.93%  "
$ws.Range("D6").Value = "'0.5203"
$ws.Range("E6").Value = "'  -0.95%  "
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("D8").Value = "'0.2565"
$ws.Range("E8").Value = "'  -3.53%  "
$ws.Range("D9").Value = "'0.06243"
$ws.Range("D10").Value = "'20.29"
$ws.Range("E10").Value = "'  -5.16%  "
$ws.Range("D11").Value = "'0.07550"
$ws.Range("E11").Value = "'  -0.09%  "
$ws.Range("D12").Value = "'1.627.40"
$ws.Range("E12").Value = "'  -3.17%  "
$ws.Range("D13").Value = "'4.356"
$ws.Range("E13").Value = "'  -2.28%  "
$ws.Range("D14").Value = "'1.859.13"
$ws.Range("D15").Value = "'0.5415"
$ws.Range("E15").Value = "'  -4.21%  "
$ws.Range("D16").Value = "'0.0₅7922"
$ws.Range("E16").Value = "'  -1.49%  "
$ws.Range("D17").Value = "'64.46"
$ws.Range("E17").Value = "'  -3.68%  "
$ws.Range("D18").Value = "'25.911.37"
$ws.Range("E18").Value = "'  -1.34%  "
$ws.Range("D20").Value = "'4.624"
$ws.Range("E20").Value = "'  -4.40%  "
$ws.Range("D21").Value = "'184.17"
$ws.Range("E21").Value = "'  -2.10%  "
$ws.Range("D22").Value = "'9.999"
$ws.Range("E22").Value = "'  -4.26%  "
$ws.Range("E23").Value = "'  -2.12%  "
$ws.Range("E24").Value = "'  -0.07%  "
$ws.Range("D25").Value = "'145.52"
$ws.Range("E25").Value = "'  -2.67%  "
$ws.Range("D26").Value = "'0.1199"
$ws.Range("E26").Value = "'  -4.21%  "
$ws.Range("D27").Value = "'7.333"
$ws.Range("E27").Value = "'  -3.45%  "
$ws.Range("D28").Value = "'15.47"
$ws.Range("E28").Value = "'  -3.56%  "
$ws.Range("E29").Value = "'  +0.96%  "
$ws.Range("D30").Value = "'0.05945"
$ws.Range("E30").Value = "'  -3.81%  "
$ws.Range("D31").Value = "'1.241"
$ws.Range("E31").Value = "'  -3.28%  "
$ws.Range("D32").Value = "'3.347"
$ws.Range("E32").Value = "'  -2.79%  "
$ws.Range("D33").Value = "'3.341"
$ws.Range("E33").Value = "'  -4.44%  "
$ws.Range("D34").Value = "'1.604"
$ws.Range("E34").Value = "'  -1.98%  "
$ws.Range("D35").Value = "'0.9699"
$ws.Range("E35").Value = "'  -3.31%  "
$ws.Range("E36").Value = "'  -0.97%  "
$ws.Range("D37").Value = "'2.729"
$ws.Range("E37").Value = "'  -0.53%  "
$ws.Range("D38").Value = "'0.5804"
$ws.Range("E38").Value = "'  -4.35%  "
$ws.Range("E39").Value = "'  -1.51%  "
$ws.Range("E40").Value = "'  -0.49%  "
$ws.Range("D41").Value = "'0.8392"
$ws.Range("E41").Value = "'  -3.40%  "
$ws.Range("D42").Value = "'1.024.94"
$ws.Range("E42").Value = "'  -5.43%  "
$ws.Range("D43").Value = "'5.656"
$ws.Range("E43").Value = "'  -7.29%  "
$ws.Range("D44").Value = "'99.63"
$ws.Range("E44").Value = "'  -0.48%  "
$ws.Range("D45").Value = "'1.784.69"
$ws.Range("E45").Value = "'  -2.40%  "
$ws.Range("D46").Value = "'0.0₈106"
$ws.Range("E46").Value = "'  -4.81%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "'  +0.25%  "
$ws.Range("D48").Value = "'54.25"
$ws.Range("E48").Value = "'  -3.70%  "
$ws.Range("D49").Value = "'7.967"
$ws.Range("E49").Value = "'  -0.45%  "
$ws.Range("D50").Value = "'0.05178"
$ws.Range("E50").Value = "'  -1.19%  "
$ws.Range("D51").Value = "'0.4225"
$ws.Range("E51").Value = "'  -0.73%  "

$touched = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D20", "E20", "D21", "E21", "D22", "E22", "E23", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "E36", "D37", "E37", "D38", "E38", "E39", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $touched) {
    $ws.Range($addr).Style = "Normal"
}
